$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.264.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "'2.013.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").Value = "'252.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.10%  "

$ws.Range("D6").Value = "'0.642"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.87%  "

$ws.Range("D7").Value = "'62.82"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.91%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "'59.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.36%  "

$ws.Range("E10").Value = "  +1.58%  "

$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("D13").Value = "'0.907"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").Value = "'14.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.57%  "

$ws.Range("D15").Value = "'2.307.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.27%  "

$ws.Range("E16").Value = "  +1.75%  "

$ws.Range("D17").Value = "'19.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +15.67%  "

$ws.Range("D18").Value = "'2.006.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("D19").Value = "'36.157.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "

$ws.Range("D20").Value = "'72.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("E21").Value = "  +1.07%  "

$ws.Range("D22").Value = "'5.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.27%  "

$ws.Range("D23").Value = "'234.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.88%  "

$ws.Range("D24").Value = "'2.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +20.76%  "

$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("E26").Value = "  -1.92%  "

$ws.Range("E27").Value = "  +4.15%  "

$ws.Range("D28").Value = "'162.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("D29").Value = "'19.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.113"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +28.98%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.120"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").Value = "'5.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.27%  "

$ws.Range("E33").Value = "  -0.65%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.17%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'2.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +16.07%  "

$ws.Range("D36").Value = "'0.0605"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.70%  "

$ws.Range("E37").Value = "  -0.37%  "

$ws.Range("D38").Value = "'1.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").Value = "'5.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.80%  "

$ws.Range("D40").Value = "'0.102"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.94%  "

$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("E42").Value = "  +1.55%  "

$ws.Range("D43").Value = "'0.0216"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("E44").Value = "  +3.25%  "

$ws.Range("E45").Value = "  +6.20%  "

$ws.Range("D46").Value = "'7.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.27%  "

$ws.Range("D47").Value = "'94.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.84%  "

$ws.Range("D48").Value = "'1.426.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.77%  "

$ws.Range("D49").Value = "'2.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.15%  "

$ws.Range("E50").Value = "  -0.22%  "

$ws.Range("D51").Value = "'47.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.48%  "
